$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial for every data row (rows 2-331).
# All of them currently hold 45171 (2023-09-02) and need to be bumped to
# 45172 (2023-09-03).
for ($r = 2; $r -le 331; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value = 45172
    }
}
